$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column K (2021) that mirrors the formatting of column J, then fill in the
# values for that year for each indicator row.

# Copy formatting (styles) from column J cells down into column K first, row by row,
# matching only the rows that already carry a style in column J (rows 3-10).
$ws.Range("J3").Copy($ws.Range("K3"))
$ws.Range("J4").Copy($ws.Range("K4"))
$ws.Range("J5").Copy($ws.Range("K5"))
$ws.Range("J6").Copy($ws.Range("K6"))
$ws.Range("J7").Copy($ws.Range("K7"))
$ws.Range("J8").Copy($ws.Range("K8"))
$ws.Range("J9").Copy($ws.Range("K9"))
$ws.Range("J10").Copy($ws.Range("K10"))

# Now set the actual cell values/content for the new 2021 column.
$ws.Range("K4").Value = 2021
$ws.Range("K5").Value = 375
$ws.Range("K6").Value = "-"
$ws.Range("K7").Value = 5
$ws.Range("K8").Value = "-"
$ws.Range("K9").Value = 18
$ws.Range("K10").Value = 150

# Match the author's last active selection in the saved file.
$ws.Range("K7").Select()
